$wb = $excel.ActiveWorkbook

# Sheet "상복" (4th sheet): update beer quantity
$ws4 = $wb.Worksheets.Item("상복")
$ws4.Range("C11").Value = 0

# Sheet "기타" (5th sheet): update quantities in column C
$ws5 = $wb.Worksheets.Item("기타")
$ws5.Range("C2").Value = 1
$ws5.Range("C3").Value = 0
$ws5.Range("C4").Value = 0
$ws5.Range("C5").Value = 2
$ws5.Range("C7").Value = 2
$ws5.Range("C8").Value = 20
$ws5.Range("C10").Value = 24
$ws5.Range("C11").Value = 2
$ws5.Range("C14").Value = 25

$wb.Save()
